$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate, in document order, the *second* occurrence of "Run #" (the table
# cell that still has the text split across two runs "Run " + "#", with a
# lastRenderedPageBreak + the hidden _GoBack bookmark sitting between them)
# and the *second* occurrence of the "Accuracy on test data" paragraph
# (the one that currently has the run directly inside it, no split yet).
# ---------------------------------------------------------------------------

function Find-AllRanges([string]$text) {
    $results = @()
    $probe = $d.Content
    $f = $probe.Find
    $f.Text = $text
    $f.Forward = $true
    $f.Wrap = 0
    $f.MatchCase = $false
    $f.MatchWholeWord = $false
    $f.MatchWildcards = $false
    while ($f.Execute()) {
        $results += ,@($probe.Start, $probe.End)
        $probe.Collapse(0)  # wdCollapseEnd
    }
    return $results
}

$runHashPositions = Find-AllRanges("Run #")
$accuracyPositions = Find-AllRanges("Accuracy on test data")

$runHashTarget = $runHashPositions[$runHashPositions.Length - 1]
$accuracyTarget = $accuracyPositions[$accuracyPositions.Length - 1]

# ---------------------------------------------------------------------------
# Change 2 (do this one first - it sits *after* the paragraph touched by
# change 1, so doing it first keeps the change-1 offsets stable): merge the
# "Run " + "#" runs in the table cell into a single run "Run #". This also
# drops the lastRenderedPageBreak and the old _GoBack bookmark that were
# sitting between the two runs, matching the diff.
# ---------------------------------------------------------------------------
$runHashRange = $d.Range($runHashTarget[0], $runHashTarget[1])
[void]$runHashRange.Find.Execute("Run #", $false, $false, $false, $false, $false, $true, 0, $false, "Run #", 1)

# ---------------------------------------------------------------------------
# Change 1: split the paragraph that holds "Accuracy on test data" into two
# paragraphs - the first (the original one) ends up empty, the second (new)
# keeps the paragraph formatting, gains a fresh _GoBack bookmark right after
# <w:pPr>, and its run gains a <w:lastRenderedPageBreak/> right before the
# text.
# ---------------------------------------------------------------------------
$accPara = $d.Range($accuracyTarget[0], $accuracyTarget[0]).Paragraphs(1)
$splitScope = $d.Range($accPara.Range.Start, $accPara.Range.End - 1)
[void]$splitScope.Find.Execute("Accuracy on test data", $false, $false, $false, $false, $false, $true, 0, $false, "^pAccuracy on test data", 2)

# The paragraph that now holds the moved text is the one right after the
# (now empty) original paragraph.
$newPara = $accPara.Next()

# Restore the run formatting (Times New Roman, 14pt / sz 28) that Find and
# Replace resets to the document default.
$fnt = $newPara.Range.Font
$fnt.Name = "Times New Roman"
$fnt.NameAscii = "Times New Roman"
$fnt.NameBi = "Times New Roman"
$fnt.Size = 14
$fnt.SizeBi = 14

# Replace the paragraph's text range with its final OOXML: a _GoBack
# bookmark followed by a single run containing the lastRenderedPageBreak and
# the text - all in one shot, so the text stays a single <w:r>.
$finalRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$finalXml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Accuracy on test data</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$finalRange.InsertXML($finalXml)

Write-Output "Done."
